$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1: Finals CSK vs KKR match result entry (row 69)
# ---------------------------------------------------------------------------
$ws1.Range("C69").Value = "CSK vs KKR"
$ws1.Range("E69").Value = 100
$ws1.Range("H69").Value = 60
$ws1.Range("K69").Value = 20
$ws1.Range("N69").Value = 40
$ws1.Range("Q69").Value = 0
$ws1.Range("T69").Value = 80

# ---------------------------------------------------------------------------
# Sheet1: fill in missing "F" column scores for the winner-prediction coin
# table (rows 80-85)
# ---------------------------------------------------------------------------
$ws1.Range("F80").Value = 3
$ws1.Range("F81").Value = 11
$ws1.Range("F82").Value = 3
$ws1.Range("F83").Value = 11
$ws1.Range("F85").Value = 11

$wb.Application.Calculate()

# ---------------------------------------------------------------------------
# Sheet2: update the Finals prediction section (rows 52-57)
# ---------------------------------------------------------------------------
$ws2.Range("F52").Value = "Finals CSK vs KKR"

$ws2.Range("G53").Value = "Sampath M"
$ws2.Range("H53").Value = "Jayanth"

$ws2.Range("G54").Value = "Anantha"
$ws2.Range("H54").Value = "Rapaka"

$ws2.Range("G55").Value = "Sampath M"
$ws2.Range("H55").Value = "Justin"

$ws2.Range("G56").Value = "Anantha"
$ws2.Range("H56").Value = "Anantha"

$ws2.Range("G57").Value = "Anantha"
$ws2.Range("H57").Value = "Rapaka"

$ws2.Range("J52").Value = "Rank 1"
$ws2.Range("K52").Value = "Anantha"
$ws2.Range("J53").Value = "Rank 2"
$ws2.Range("K53").Value = "Sampath M"

# ---------------------------------------------------------------------------
# Sheet2: new final standings table (J60:M66)
# ---------------------------------------------------------------------------
$ws2.Range("J60").Value = "Rank"
$ws2.Range("K60").Value = "Player"
$ws2.Range("L60").Value = "Total"

$ws2.Range("J61").Value = 1
$ws2.Range("K61").Value = "Anantha"
$ws2.Range("L61").Value = 328.75
$ws2.Range("M61").Value = "Congrats"

$ws2.Range("J62").Value = 2
$ws2.Range("K62").Value = "Sushma"
$ws2.Range("L62").Value = 85
$ws2.Range("M62").Value = "Congrats"

$ws2.Range("J63").Value = 3
$ws2.Range("K63").Value = "Rapaka"
$ws2.Range("L63").Value = 15.3125
$ws2.Range("M63").Value = "Congrats"

$ws2.Range("J64").Value = 4
$ws2.Range("K64").Value = "Sampath M"
$ws2.Range("L64").Value = -74.0625

$ws2.Range("J65").Value = 5
$ws2.Range("K65").Value = "Jayanth"
$ws2.Range("L65").Value = -117.5

$ws2.Range("J66").Value = 6
$ws2.Range("K66").Value = "Justin"
$ws2.Range("L66").Value = -237.5

# Formatting: copy matching existing styles onto the new cells so the
# workbook's style table grows the same way Excel would grow it.
$ws1.Range("A10").Copy()
$ws2.Range("J61:J66").PasteSpecial(-4122)

$ws1.Range("L80").Copy()
$ws2.Range("K61:K66").PasteSpecial(-4122)
$ws2.Range("R80").Copy()
$ws2.Range("J52").PasteSpecial(-4122)
$ws2.Range("J53").PasteSpecial(-4122)

$ws1.Range("M80").Copy()
$ws2.Range("L61:L66").PasteSpecial(-4122)

$ws2.Range("H38").Copy()
$ws2.Range("J60:L60").PasteSpecial(-4122)
$ws2.Range("J60:L60").Font.Bold = $true
$ws2.Range("M61:M63").Value = $ws2.Range("M61:M63").Value
$ws2.Range("J60").Copy()
$ws2.Range("M61:M63").PasteSpecial(-4122)
$ws2.Range("M61").Value = "Congrats"
$ws2.Range("M62").Value = "Congrats"
$ws2.Range("M63").Value = "Congrats"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Conditional formatting on the new standings column, matching the style
# already used for the similar coin-total columns on Sheet1.
# ---------------------------------------------------------------------------
$ws1.Range("Q85").FormatConditions.Copy($ws2.Range("L66"))
$ws1.Range("M80:M84").FormatConditions.Copy($ws2.Range("L61:L65"))

# ---------------------------------------------------------------------------
# View state: Sheet2 becomes the active/selected sheet, Sheet1 keeps its
# frozen-pane selection pointed at the coin-prediction table.
# ---------------------------------------------------------------------------
$ws1.Range("Q80:R85").Select()
$ws2.Activate()
$ws2.Range("P67").Select()

$wb.Application.Calculate()
